$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The timetable already has a Monday (rows 2-8) and Tuesday (rows 9-15) block
# filled in for columns C:F. Replicate the Tuesday block's class entries into
# the Wednesday (16-22), Thursday (23-29) and Friday (30-36) blocks.

$startRows = @(16, 23, 30)

foreach ($startRow in $startRows) {
    $ws.Cells.Item($startRow, 3).Value = "PL-204"
    $ws.Cells.Item($startRow, 5).Value = "AM-BT"
    $ws.Cells.Item($startRow, 6).Value = "Fizica-192"

    $ws.Cells.Item($startRow + 1, 3).Value = "PC-205"

    $ws.Cells.Item($startRow + 2, 4).Value = "ENG-BT"
    $ws.Cells.Item($startRow + 2, 5).Value = "PL-204"

    $ws.Cells.Item($startRow + 3, 5).Value = "PC-205"
    $ws.Cells.Item($startRow + 3, 6).Value = "ENG-BT"

    $ws.Cells.Item($startRow + 4, 4).Value = "PL-204"

    $ws.Cells.Item($startRow + 5, 4).Value = "PC-205"
}

# Move the active selection as recorded in the saved session state.
$ws.Range("I26").Select()
